# Auto-generated Excel COM-interop script applying the data refresh
# described in the commit "Update data: 2025-10-29 18:28"

$wb = $excel.ActiveWorkbook

# ---- "Top Gainers" sheet: update rows 42-76 (Stock / Col1 / Col2 / Col3) ----
$wsGainers = $wb.Worksheets.Item("Top Gainers")

$wsGainers.Range("B42").Value = "INDOTHAI"
$wsGainers.Range("C42").Value = 4.8064
$wsGainers.Range("D42").Value = 4.5349
$wsGainers.Range("E42").Value = 43.748

$wsGainers.Range("B43").Value = "SANDUMA"
$wsGainers.Range("C43").Value = 4.593
$wsGainers.Range("D43").Value = 2.1405
$wsGainers.Range("E43").Value = 30.2813

$wsGainers.Range("B44").Value = "LLOYDSENT"
$wsGainers.Range("C44").Value = 4.5646
$wsGainers.Range("D44").Value = 1.8339
$wsGainers.Range("E44").Value = 11.234

$wsGainers.Range("B45").Value = "STAR"
$wsGainers.Range("C45").Value = 4.5025
$wsGainers.Range("D45").Value = 4.4319
$wsGainers.Range("E45").Value = 3.662

$wsGainers.Range("B46").Value = "RECLTD"
$wsGainers.Range("C46").Value = 4.4992
$wsGainers.Range("D46").Value = 3.4756
$wsGainers.Range("E46").Value = 3.4062

$wsGainers.Range("B47").Value = "NBCC"
$wsGainers.Range("C47").Value = 4.4511
$wsGainers.Range("D47").Value = 3.1605
$wsGainers.Range("E47").Value = 7.6018

$wsGainers.Range("B48").Value = "GPPL"
$wsGainers.Range("C48").Value = 4.4154
$wsGainers.Range("D48").Value = 3.4073
$wsGainers.Range("E48").Value = 5.0497

$wsGainers.Range("B49").Value = "BIL"
$wsGainers.Range("C49").Value = 4.3654
$wsGainers.Range("D49").Value = 9.122199999999999
$wsGainers.Range("E49").Value = -0.3203

$wsGainers.Range("B50").Value = "HUDCO"
$wsGainers.Range("C50").Value = 4.3201
$wsGainers.Range("D50").Value = 3.8924
$wsGainers.Range("E50").Value = 5.3884

$wsGainers.Range("B51").Value = "SGMART"
$wsGainers.Range("C51").Value = 4.2736
$wsGainers.Range("D51").Value = 8.258900000000001
$wsGainers.Range("E51").Value = 2.5381

$wsGainers.Range("B52").Value = "MRPL"
$wsGainers.Range("C52").Value = 4.2642
$wsGainers.Range("D52").Value = 9.7103
$wsGainers.Range("E52").Value = 20.0542

$wsGainers.Range("B53").Value = "JKIL"
$wsGainers.Range("C53").Value = 4.1372
$wsGainers.Range("D53").Value = 2.9463
$wsGainers.Range("E53").Value = 1.7584

$wsGainers.Range("B54").Value = "SAMBHV"
$wsGainers.Range("C54").Value = 4.1349
$wsGainers.Range("D54").Value = 2.624
$wsGainers.Range("E54").Value = 5.167

$wsGainers.Range("B55").Value = "SAPPHIRE"
$wsGainers.Range("C55").Value = 4.1265
$wsGainers.Range("D55").Value = 1.7633
$wsGainers.Range("E55").Value = -0.7999000000000001

$wsGainers.Range("B56").Value = "PVRINOX"
$wsGainers.Range("C56").Value = 4.1118
$wsGainers.Range("D56").Value = 6.2102
$wsGainers.Range("E56").Value = 14.707

$wsGainers.Range("B57").Value = "KERNEX"
$wsGainers.Range("C57").Value = 3.9981
$wsGainers.Range("D57").Value = 7.4592
$wsGainers.Range("E57").Value = 27.1054

$wsGainers.Range("B58").Value = "SUNFLAG"
$wsGainers.Range("C58").Value = 3.997
$wsGainers.Range("D58").Value = 4.333
$wsGainers.Range("E58").Value = 4.6312

$wsGainers.Range("B59").Value = "CMSINFO"
$wsGainers.Range("C59").Value = 3.9096
$wsGainers.Range("D59").Value = 2.6872
$wsGainers.Range("E59").Value = 2.8935

$wsGainers.Range("B60").Value = "GMBREW"
$wsGainers.Range("C60").Value = 3.8999
$wsGainers.Range("D60").Value = -0.53
$wsGainers.Range("E60").Value = 79.029

$wsGainers.Range("B61").Value = "APARINDS"
$wsGainers.Range("C61").Value = 3.8924
$wsGainers.Range("D61").Value = 8.3414
$wsGainers.Range("E61").Value = 15.5876

$wsGainers.Range("B62").Value = "HITECHGEAR"
$wsGainers.Range("C62").Value = 3.8587
$wsGainers.Range("D62").Value = 1.1486
$wsGainers.Range("E62").Value = 9.9254

$wsGainers.Range("B63").Value = "NPST"
$wsGainers.Range("C63").Value = 3.8509
$wsGainers.Range("D63").Value = -2.0059
$wsGainers.Range("E63").Value = -3.5057

$wsGainers.Range("B67").Value = "DCW"
$wsGainers.Range("C67").Value = 3.7544
$wsGainers.Range("D67").Value = 2.3219
$wsGainers.Range("E67").Value = -3.9753

$wsGainers.Range("B68").Value = "RHETAN"
$wsGainers.Range("C68").Value = 3.754
$wsGainers.Range("D68").Value = 4.178
$wsGainers.Range("E68").Value = 6.549

$wsGainers.Range("B69").Value = "HINDPETRO"
$wsGainers.Range("C69").Value = 3.6935
$wsGainers.Range("D69").Value = 6.9335
$wsGainers.Range("E69").Value = 5.7397

$wsGainers.Range("B70").Value = "BHARTIHEXA"
$wsGainers.Range("C70").Value = 3.6718
$wsGainers.Range("D70").Value = 7.0877
$wsGainers.Range("E70").Value = 15.3332

$wsGainers.Range("B71").Value = "HLEGLAS"
$wsGainers.Range("C71").Value = 3.659
$wsGainers.Range("D71").Value = 8.115500000000001
$wsGainers.Range("E71").Value = 27.1239

$wsGainers.Range("B72").Value = "RHIM"
$wsGainers.Range("C72").Value = 3.6544
$wsGainers.Range("D72").Value = 3.2276
$wsGainers.Range("E72").Value = 5.1826

$wsGainers.Range("B73").Value = "SHK"
$wsGainers.Range("C73").Value = 3.6347
$wsGainers.Range("D73").Value = 2.388
$wsGainers.Range("E73").Value = -1.932

$wsGainers.Range("B74").Value = "BCLIND"
$wsGainers.Range("C74").Value = 3.6271
$wsGainers.Range("D74").Value = 2.2945
$wsGainers.Range("E74").Value = 0.1728

$wsGainers.Range("B75").Value = "MUKANDLTD"
$wsGainers.Range("C75").Value = 3.6133
$wsGainers.Range("D75").Value = 11.9685
$wsGainers.Range("E75").Value = 9.550800000000001

$wsGainers.Range("B76").Value = "CGPOWER"
$wsGainers.Range("C76").Value = 3.6125
$wsGainers.Range("D76").Value = 3.4192
$wsGainers.Range("E76").Value = 1.0325

# ---- "1 Month Performance" sheet: update rows 16-19, 36, 72-73 (Stock / % Change) ----
$wsPerf = $wb.Worksheets.Item("1 Month Performance")

$wsPerf.Range("B16").Value = "V2RETAIL"
$wsPerf.Range("C16").Value = 37.2004

$wsPerf.Range("B17").Value = "RAMAPHO"
$wsPerf.Range("C17").Value = 36.9731

$wsPerf.Range("B18").Value = "SANDUMA"
$wsPerf.Range("C18").Value = 36.9057

$wsPerf.Range("B19").Value = "SEJALLTD"
$wsPerf.Range("C19").Value = 36.8123

$wsPerf.Range("C36").Value = 27.4033

$wsPerf.Range("B72").Value = "KARURVYSYA"
$wsPerf.Range("C72").Value = 19.11

$wsPerf.Range("B73").Value = "IIFL"
$wsPerf.Range("C73").Value = 18.9853
